$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add row 43 ("tp 7 and 8" from OR tide pools) mirroring the layout of row 42.
# Copy formatting from A42 so the date cell picks up the existing date style
# instead of minting a new numFmt/style entry.
$ws.Range("A42").Copy()
$ws.Range("A43").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("A43").Value = 43713
$ws.Range("B43").Value = 2219.0529999999999
$ws.Range("C43").Value = 2207.0300000000002
$ws.Range("D43").Formula = "=100*(B43-C43)/C43"
$ws.Range("E43").Value = 169
$ws.Range("F43").Value = "Opened CRM (8/30/2019)"

# Reflect the new last row in the sheet view (scroll + selection).
$ws.Application.ActiveWindow.ScrollRow = 32
$ws.Range("F42:F43").Select()
